$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-3: name / card number updates
$ws.Range("C2").Value = "Hartmut"

# B3 holds a purely-numeric-looking card number that must stay text (like
# the original). Force text storage via NumberFormat, then restore the
# original cell style (lost when the number format changed) by pasting
# just the formatting back in from an unmodified sibling cell that still
# carries the same style.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("C3").Value = "Mohaupt"

# Row 5: starting balance date
$ws.Range("D5").Value = "KONTOSTAND AM 04.07.2024"

# Row 6
$ws.Range("B6").Value = "05.07."
$ws.Range("C6").Value = "06.07."
$ws.Range("D6").Value = "KARTENZAHLUNG JET TANKSTELLE"
$ws.Range("E6").Value = "51,39-"

# Row 7
$ws.Range("B7").Value = "06.07."
$ws.Range("C7").Value = "07.07."
$ws.Range("D7").Value = "MCDONALDS Soltau"
$ws.Range("E7").Value = "44,17-"

# Row 8
$ws.Range("B8").Value = "07.07."
$ws.Range("C8").Value = "08.07."
$ws.Range("D8").Value = "PAYPAL IMSJGM"
$ws.Range("E8").Value = "48,42-"

# Row 9 (newly populated, style for E9 must match E6/E7/E8 -> copy style)
$ws.Range("B9").Value = "11.07."
$ws.Range("C9").Value = "12.07."
$ws.Range("D9").Value = "RECHNUNG VODAFONE GMBH 12030635"
$ws.Range("E9").Value = "41,16-"
$ws.Range("E8").Copy()
$ws.Range("E9").PasteSpecial(-4122)  # xlPasteFormats

# Row 10 (newly populated, style for E10 must match E6/E7/E8 -> copy style)
$ws.Range("B10").Value = "15.07."
$ws.Range("C10").Value = "16.07."
$ws.Range("D10").Value = "BURGER KING Neustrelitz"
$ws.Range("E10").Value = "29,23-"
$ws.Range("E8").Copy()
$ws.Range("E10").PasteSpecial(-4122)  # xlPasteFormats

# Row 12: ending balance date + amount
$ws.Range("D12").Value = "KONTOSTAND AM 19.07.2024"
$ws.Range("E12").Value = "214,37-"

# Row 13: next billing date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 27.07.2024"

$excel.CutCopyMode = $false
